# Adding new RAAL model Production
# Updates computed GHI production values on the "Daily" and "Hourly" sheets.

$wb = $excel.ActiveWorkbook

# ---- Daily sheet ----
$daily = $wb.Worksheets.Item("Daily")

$daily.Range("G2").Value = 2848.79
$daily.Range("H2").Value = 6099.75
$daily.Range("I2").Value = 712.34
$daily.Range("J2").Value = 768.91
$daily.Range("K2").Value = 0
$daily.Range("L2").Value = 768.91

# ---- Hourly sheet ----
$hourly = $wb.Worksheets.Item("Hourly")

# Row 9
$hourly.Range("I9").Value = 22.27
$hourly.Range("K9").Value = 0.78
$hourly.Range("M9").Value = 0.78

# Row 10
$hourly.Range("H10").Value = 86.78
$hourly.Range("I10").Value = 383.2
$hourly.Range("K10").Value = 24.48
$hourly.Range("L10").Value = 0
$hourly.Range("M10").Value = 24.48

# Row 11
$hourly.Range("H11").Value = 228.42
$hourly.Range("I11").Value = 615.79
$hourly.Range("J11").Value = 69.45999999999999
$hourly.Range("K11").Value = 61.64
$hourly.Range("L11").Value = 0
$hourly.Range("M11").Value = 61.64

# Row 12
$hourly.Range("H12").Value = 352.66
$hourly.Range("I12").Value = 722.74
$hourly.Range("J12").Value = 83.97
$hourly.Range("K12").Value = 102.45
$hourly.Range("L12").Value = 0
$hourly.Range("M12").Value = 102.45

# Row 13
$hourly.Range("H13").Value = 438
$hourly.Range("I13").Value = 775.6799999999999
$hourly.Range("J13").Value = 92.03
$hourly.Range("K13").Value = 121.94
$hourly.Range("L13").Value = 0
$hourly.Range("M13").Value = 121.94

# Row 14
$hourly.Range("H14").Value = 473.78
$hourly.Range("I14").Value = 794.74
$hourly.Range("J14").Value = 95.08
$hourly.Range("K14").Value = 128.88
$hourly.Range("L14").Value = 0
$hourly.Range("M14").Value = 128.88

# Row 15
$hourly.Range("H15").Value = 455.95
$hourly.Range("I15").Value = 785.53
$hourly.Range("J15").Value = 93.56
$hourly.Range("K15").Value = 122.13
$hourly.Range("L15").Value = 0
$hourly.Range("M15").Value = 122.13

# Row 16
$hourly.Range("H16").Value = 386.5
$hourly.Range("I16").Value = 745.39
$hourly.Range("J16").Value = 87.27
$hourly.Range("K16").Value = 100.55
$hourly.Range("M16").Value = 100.55

# Row 17
$hourly.Range("H17").Value = 273.73
$hourly.Range("I17").Value = 660.6799999999999
$hourly.Range("J17").Value = 75.2
$hourly.Range("K17").Value = 68.43000000000001
$hourly.Range("M17").Value = 68.43000000000001

# Row 18
$hourly.Range("H18").Value = 134.28
$hourly.Range("I18").Value = 486.5
$hourly.Range("K18").Value = 33.57
$hourly.Range("M18").Value = 33.57

# Row 19
$hourly.Range("I19").Value = 107.24
$hourly.Range("K19").Value = 4.05
$hourly.Range("M19").Value = 4.05
